$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 11, column A: "LIU " -> "LIU" (trailing space removed)
$ws.Range("A11").Value = "LIU"

# New header cells in row 1 (use the same bold style as the existing header cells)
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Role"
$ws.Range("F1").Value = "Committee Of"
$ws.Range("D1:F1").Font.Bold = $true

# Widen the new "Committee Of" column
$ws.Columns.Item(6).ColumnWidth = 15.65

# Fill in Password / Role columns for each student row (2-12)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 4).Value = "password"
    $ws.Cells.Item($r, 5).Value = "Student"
}

# Row 11 (LIU) has a special password value of "1" stored as text
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1"
$ws.Range("D11").Style = "Normal"

# Update selection to match target workbook view state
$ws.Range("F1:F1048576").Select()
